$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 244, shifting existing rows 244:297 down to 245:298
$ws.Rows("244:244").Insert()

# Populate the newly inserted row 244 with the new data record
$ws.Range("A244").Value = 8
$ws.Range("B244").Value = "Terminal La Palmera de La Serena"
$ws.Range("C244").Value = "Coquimbo"
$ws.Range("D244").Value = 44543
$ws.Range("D244").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E244").Value = 4
$ws.Range("F244").Value = 100114001
$ws.Range("G244").Value = "Papa"
$ws.Range("H244").Value = "Asterix"
$ws.Range("I244").Value = "1a nueva(o)"
$ws.Range("J244").Value = 2000
$ws.Range("K244").Value = 10000
$ws.Range("L244").Value = 11000
$ws.Range("M244").Value = 10500
$ws.Range("N244").Value = "$/saco 25 kilos"
$ws.Range("O244").Value = "Provincia de Melipilla"
$ws.Range("P244").Value = 420
$ws.Range("Q244").Value = 25
$ws.Range("R244").Value = "Hortaliza"
